$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")
$ws.Activate()

# Fill in boolean TRUE ("Test Result" = passed) for every customer row (2-14),
# simulating a CRUD test run that validated each customer record.
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 6).Value = $true
}

# Scroll the view down a bit (row 7 becomes the top visible row) and leave
# the selection on a single cell, F2, as it was left after the test run.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F2").Select()
